$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: clear F2
$ws.Range("F2").Value = ""

# Row 6: set K6
$ws.Range("K6").Value = "7,97 TL - 15,96 TL - 199,41 TL"

# Row 7: clear F7
$ws.Range("F7").Value = ""

# Row 12: update D12, set K12
$ws.Range("D12").Value = "WU: 0,75 USD–12 USD; Diğer: 700 TL–4.000 TL"
$ws.Range("K12").Value = "WU: ,USD–; Diğer: 529 TL–4.454,74 TL"

# Row 13: update C13, set K13
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

# Row 14: set K14
$ws.Range("K14").Value = "1.196,51 TL - 5.583,74 TL"
